$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 20
$ws.Range("I8").Value = 20
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 60
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 79
$ws.Range("N8").Value = $null
$ws.Range("H17").Value = 2073.4583
$ws.Range("J17").Value = 2095.3809
$ws.Range("L17").Value = 6286.1427
$ws.Range("N17").Value = -6622.1427
$ws.Range("H43").Value = 2802.7778
$ws.Range("I43").Value = 3207.1667
$ws.Range("J43").Value = 2600.5833
$ws.Range("K43").Value = 3207.1667
$ws.Range("L43").Value = 2600.5833
$ws.Range("M43").Value = -3138.1667
$ws.Range("N43").Value = -2738.5833
$ws.Range("H74").Value = 20000.5
$ws.Range("I74").Value = 20001
$ws.Range("K74").Value = 20001
$ws.Range("M74").Value = -19065
$ws.Range("H77").Value = 20000.5
$ws.Range("I77").Value = 20001
$ws.Range("K77").Value = 100005
$ws.Range("M77").Value = -95325
$ws.Range("H80").Value = 501.2
$ws.Range("I80").Value = 501
$ws.Range("J80").Value = 501.5
$ws.Range("K80").Value = 1503
$ws.Range("L80").Value = 1504.5
$ws.Range("M80").Value = -505
$ws.Range("N80").Value = -3500.5
$ws.Range("H83").Value = 501.2
$ws.Range("I83").Value = 501
$ws.Range("J83").Value = 501.5
$ws.Range("K83").Value = 4509
$ws.Range("L83").Value = 4513.5
$ws.Range("M83").Value = 483
$ws.Range("N83").Value = -14497.5
$ws.Range("H98").Value = 3229.4783
$ws.Range("I98").Value = 3126.318
$ws.Range("K98").Value = 3126.318
$ws.Range("M98").Value = -1628.318
$ws.Range("H122").Value = 3229.4783
$ws.Range("I122").Value = 3126.318
$ws.Range("K122").Value = 9378.954000000002
$ws.Range("M122").Value = -6928.954000000002
$ws.Range("H137").Value = 4202.35
$ws.Range("I137").Value = 1617.0344
$ws.Range("K137").Value = 4851.1032
$ws.Range("M137").Value = -2301.1032
$ws.Range("H141").Value = 4547.3335
$ws.Range("J141").Value = 5399.75
$ws.Range("L141").Value = 16199.25
$ws.Range("N141").Value = -26559.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4556.209
$ws.Range("I32").Value = 3847.9375
$ws.Range("K32").Value = 3847.9375
$ws.Range("M32").Value = -3560.9375
$ws.Range("H45").Value = 45426.9
$ws.Range("I45").Value = 72640.164
$ws.Range("K45").Value = 72640.164
$ws.Range("M45").Value = -72263.164
$ws.Range("H61").Value = 4879.1787
$ws.Range("I61").Value = 1550.1111
$ws.Range("K61").Value = 1550.1111
$ws.Range("M61").Value = -1338.1111
$ws.Range("H74").Value = 187379.56
$ws.Range("I74").Value = 254160.64
$ws.Range("K74").Value = 254160.64
$ws.Range("M74").Value = -253286.64
$ws.Range("H77").Value = 187379.56
$ws.Range("I77").Value = 254160.64
$ws.Range("K77").Value = 1270803.2
$ws.Range("M77").Value = -1266435.2
$ws.Range("H88").Value = 5999.25
$ws.Range("I88").Value = 5749.5
$ws.Range("J88").Value = 6249
$ws.Range("K88").Value = 5749.5
$ws.Range("L88").Value = 6249
$ws.Range("M88").Value = -5343.5
$ws.Range("N88").Value = -7061
$ws.Range("H91").Value = 5999.25
$ws.Range("I91").Value = 5749.5
$ws.Range("J91").Value = 6249
$ws.Range("K91").Value = 5749.5
$ws.Range("L91").Value = 6249
$ws.Range("M91").Value = -4345.5
$ws.Range("N91").Value = -9057
$ws.Range("H97").Value = 1103.6666
$ws.Range("I97").Value = 1117.95
$ws.Range("J97").Value = 818
$ws.Range("K97").Value = 1117.95
$ws.Range("L97").Value = 818
$ws.Range("M97").Value = -621.95
$ws.Range("N97").Value = -1810
$ws.Range("H110").Value = 2556.238
$ws.Range("I110").Value = 1317.0667
$ws.Range("K110").Value = 1317.0667
$ws.Range("M110").Value = 727.9332999999999
$ws.Range("H122").Value = 3722.5386
$ws.Range("I122").Value = 3489.0454
$ws.Range("K122").Value = 10467.1362
$ws.Range("M122").Value = -8017.136200000001
$ws.Range("H132").Value = 2518.8857
$ws.Range("I132").Value = 1746.7037
$ws.Range("K132").Value = 5240.1111
$ws.Range("M132").Value = -2710.1111
$ws.Range("H136").Value = 4879.1787
$ws.Range("I136").Value = 1550.1111
$ws.Range("K136").Value = 4650.3333
$ws.Range("M136").Value = -2100.3333
$ws.Range("H139").Value = 81416.07000000001
$ws.Range("J139").Value = 81416.07000000001
$ws.Range("L139").Value = 81416.07000000001
$ws.Range("N139").Value = -91696.07000000001
$ws.Range("H140").Value = 87220.78
$ws.Range("J140").Value = 87220.78
$ws.Range("L140").Value = 87220.78
$ws.Range("N140").Value = -97580.78
$ws.Range("H141").Value = 135665
$ws.Range("J141").Value = 135665
$ws.Range("L141").Value = 135665
$ws.Range("N141").Value = -146025

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4101.12
$ws.Range("I86").Value = 3811.45
$ws.Range("J86").Value = 5259.8
$ws.Range("K86").Value = 3811.45
$ws.Range("L86").Value = 5259.8
$ws.Range("M86").Value = -2688.45
$ws.Range("N86").Value = -7505.8
$ws.Range("H89").Value = 4101.12
$ws.Range("I89").Value = 3811.45
$ws.Range("J89").Value = 5259.8
$ws.Range("K89").Value = 19057.25
$ws.Range("L89").Value = 26299
$ws.Range("M89").Value = -13441.25
$ws.Range("N89").Value = -37531
$ws.Range("H94").Value = 83333660
$ws.Range("I94").Value = 83333660
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 83333660
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -83333209
$ws.Range("N94").Value = $null
$ws.Range("H105").Value = 7224809.5
$ws.Range("I105").Value = 372855.53
$ws.Range("J105").Value = 27780670
$ws.Range("K105").Value = 372855.53
$ws.Range("L105").Value = 27780670
$ws.Range("M105").Value = -371108.53
$ws.Range("N105").Value = -27784164
$ws.Range("H108").Value = 74997.664
$ws.Range("I108").Value = 74997.664
$ws.Range("K108").Value = 74997.664
$ws.Range("M108").Value = -71157.664
$ws.Range("H133").Value = 67142.71000000001
$ws.Range("J133").Value = 94999.75
$ws.Range("L133").Value = 94999.75
$ws.Range("N133").Value = -105119.75
$ws.Range("H134").Value = 3779.25
$ws.Range("I134").Value = 3586.2
$ws.Range("K134").Value = 10758.6
$ws.Range("M134").Value = -8223.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4404.132
$ws.Range("I31").Value = 3361.8518
$ws.Range("J31").Value = 5486.5
$ws.Range("K31").Value = 3361.8518
$ws.Range("L31").Value = 5486.5
$ws.Range("M31").Value = -3066.8518
$ws.Range("N31").Value = -6076.5
$ws.Range("H34").Value = 4404.132
$ws.Range("I34").Value = 3361.8518
$ws.Range("J34").Value = 5486.5
$ws.Range("K34").Value = 3361.8518
$ws.Range("L34").Value = 5486.5
$ws.Range("M34").Value = -3159.8518
$ws.Range("N34").Value = -5890.5
$ws.Range("H76").Value = 8481
$ws.Range("I76").Value = 8481
$ws.Range("K76").Value = 8481
$ws.Range("M76").Value = -8166
$ws.Range("H79").Value = 8481
$ws.Range("I79").Value = 8481
$ws.Range("K79").Value = 8481
$ws.Range("M79").Value = -7389
$ws.Range("H107").Value = 615.3077
$ws.Range("I107").Value = 678.58826
$ws.Range("K107").Value = 678.58826
$ws.Range("M107").Value = 1241.41174
$ws.Range("H122").Value = 3616.6875
$ws.Range("I122").Value = 3109.4167
$ws.Range("J122").Value = 5138.5
$ws.Range("K122").Value = 9328.250100000001
$ws.Range("L122").Value = 15415.5
$ws.Range("M122").Value = -6878.250100000001
$ws.Range("N122").Value = -20315.5
$ws.Range("H141").Value = 90323.5
$ws.Range("J141").Value = 90323.5
$ws.Range("L141").Value = 90323.5
$ws.Range("N141").Value = -100683.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 6625.125
$ws.Range("I22").Value = 500.5
$ws.Range("J22").Value = 8666.666999999999
$ws.Range("K22").Value = 1501.5
$ws.Range("L22").Value = 26000.001
$ws.Range("M22").Value = -1332.5
$ws.Range("N22").Value = -26338.001
$ws.Range("H27").Value = 6625.125
$ws.Range("I27").Value = 500.5
$ws.Range("J27").Value = 8666.666999999999
$ws.Range("K27").Value = 1501.5
$ws.Range("L27").Value = 26000.001
$ws.Range("M27").Value = -1399.5
$ws.Range("N27").Value = -26204.001
$ws.Range("H51").Value = 3166.3333
$ws.Range("I51").Value = 2500
$ws.Range("K51").Value = 7500
$ws.Range("M51").Value = -7040
$ws.Range("H107").Value = 683.35297
$ws.Range("I107").Value = 505.5
$ws.Range("K107").Value = 1516.5
$ws.Range("M107").Value = 403.5
$ws.Range("H121").Value = 125280.125
$ws.Range("I121").Value = 182
$ws.Range("J121").Value = 200339
$ws.Range("K121").Value = 546
$ws.Range("L121").Value = 601017
$ws.Range("M121").Value = 764
$ws.Range("N121").Value = -603637
$ws.Range("H122").Value = 1090.4445
$ws.Range("I122").Value = 890
$ws.Range("J122").Value = 1147.7142
$ws.Range("K122").Value = 8010
$ws.Range("L122").Value = 10329.4278
$ws.Range("M122").Value = -5560
$ws.Range("N122").Value = -15229.4278
$ws.Range("H132").Value = 3075.0938
$ws.Range("I132").Value = 1916.75
$ws.Range("J132").Value = 3770.1
$ws.Range("K132").Value = 17250.75
$ws.Range("L132").Value = 33930.9
$ws.Range("M132").Value = -14720.75
$ws.Range("N132").Value = -38990.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 3748.75
$ws.Range("J36").Value = 3748.75
$ws.Range("L36").Value = 3748.75
$ws.Range("N36").Value = -4718.75
$ws.Range("H70").Value = 15220070
$ws.Range("I70").Value = 21830796
$ws.Range("J70").Value = 15397.9
$ws.Range("K70").Value = 21830796
$ws.Range("L70").Value = 15397.9
$ws.Range("M70").Value = -21830526
$ws.Range("N70").Value = -15937.9
$ws.Range("H73").Value = 15220070
$ws.Range("I73").Value = 21830796
$ws.Range("J73").Value = 15397.9
$ws.Range("K73").Value = 21830796
$ws.Range("L73").Value = 15397.9
$ws.Range("M73").Value = -21829860
$ws.Range("N73").Value = -17269.9
$ws.Range("H80").Value = 52633690
$ws.Range("I80").Value = 90910904
$ws.Range("J80").Value = 2518.25
$ws.Range("K80").Value = 90910904
$ws.Range("L80").Value = 2518.25
$ws.Range("M80").Value = -90909906
$ws.Range("N80").Value = -4514.25
$ws.Range("H83").Value = 52633690
$ws.Range("I83").Value = 90910904
$ws.Range("J83").Value = 2518.25
$ws.Range("K83").Value = 454554520
$ws.Range("L83").Value = 12591.25
$ws.Range("M83").Value = -454549528
$ws.Range("N83").Value = -22575.25
$ws.Range("H97").Value = 10335.111
$ws.Range("I97").Value = 9000.714
$ws.Range("J97").Value = 15005.5
$ws.Range("K97").Value = 9000.714
$ws.Range("L97").Value = 15005.5
$ws.Range("M97").Value = -8504.714
$ws.Range("N97").Value = -15997.5
$ws.Range("H102").Value = 1301.5
$ws.Range("I102").Value = 1126.7084
$ws.Range("K102").Value = 1126.7084
$ws.Range("M102").Value = 495.2916
$ws.Range("H122").Value = 2386.182
$ws.Range("I122").Value = 2414.55
$ws.Range("K122").Value = 7243.650000000001
$ws.Range("M122").Value = -4793.650000000001
$ws.Range("H126").Value = 3290.0588
$ws.Range("I126").Value = 2085.5
$ws.Range("J126").Value = 6181
$ws.Range("K126").Value = 6256.5
$ws.Range("L126").Value = 18543
$ws.Range("M126").Value = -3786.5
$ws.Range("N126").Value = -23483
$ws.Range("H132").Value = 2849.2
$ws.Range("I132").Value = 2410.037
$ws.Range("J132").Value = 3761.3076
$ws.Range("K132").Value = 7230.110999999999
$ws.Range("L132").Value = 11283.9228
$ws.Range("M132").Value = -4700.110999999999
$ws.Range("N132").Value = -16343.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 587.8
$ws.Range("I22").Value = 599.875
$ws.Range("J22").Value = 574
$ws.Range("K22").Value = 599.875
$ws.Range("L22").Value = 574
$ws.Range("M22").Value = -304.875
$ws.Range("N22").Value = -1164
$ws.Range("H27").Value = 587.8
$ws.Range("I27").Value = 599.875
$ws.Range("J27").Value = 574
$ws.Range("K27").Value = 599.875
$ws.Range("L27").Value = 574
$ws.Range("M27").Value = -492.875
$ws.Range("N27").Value = -788
$ws.Range("H40").Value = 5529.5
$ws.Range("I40").Value = 5364.1113
$ws.Range("J40").Value = 9995
$ws.Range("K40").Value = 5364.1113
$ws.Range("L40").Value = 9995
$ws.Range("M40").Value = -5228.1113
$ws.Range("N40").Value = -10267
$ws.Range("H61").Value = 2016.9286
$ws.Range("I61").Value = 2107.875
$ws.Range("K61").Value = 2107.875
$ws.Range("M61").Value = -1905.875
$ws.Range("H82").Value = 2438.7334
$ws.Range("I82").Value = 2987.1667
$ws.Range("J82").Value = 245
$ws.Range("K82").Value = 2987.1667
$ws.Range("L82").Value = 245
$ws.Range("M82").Value = -2626.1667
$ws.Range("N82").Value = -967
$ws.Range("H85").Value = 2438.7334
$ws.Range("I85").Value = 2987.1667
$ws.Range("J85").Value = 245
$ws.Range("K85").Value = 2987.1667
$ws.Range("L85").Value = 245
$ws.Range("M85").Value = -1739.1667
$ws.Range("N85").Value = -2741
$ws.Range("H113").Value = 2016.9286
$ws.Range("I113").Value = 2107.875
$ws.Range("K113").Value = 2107.875
$ws.Range("M113").Value = 62.125
$ws.Range("H122").Value = 4004.8518
$ws.Range("I122").Value = 3946.889
$ws.Range("J122").Value = 4120.778
$ws.Range("K122").Value = 11840.667
$ws.Range("L122").Value = 12362.334
$ws.Range("M122").Value = -9390.667000000001
$ws.Range("N122").Value = -17262.334
$ws.Range("H138").Value = 69405.8
$ws.Range("J138").Value = 69405.8
$ws.Range("L138").Value = 69405.8
$ws.Range("N138").Value = -79685.8
$ws.Range("H141").Value = 93299.664
$ws.Range("J141").Value = 93299.664
$ws.Range("L141").Value = 93299.664
$ws.Range("N141").Value = -103659.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 64111.5
$ws.Range("J46").Value = 64111.5
$ws.Range("L46").Value = 64111.5
$ws.Range("N46").Value = -64573.5
$ws.Range("H122").Value = 8623943
$ws.Range("I122").Value = 3034.375
$ws.Range("K122").Value = 9103.125
$ws.Range("M122").Value = -6653.125
$ws.Range("H126").Value = 1611.8572
$ws.Range("I126").Value = 1611.8572
$ws.Range("K126").Value = 4835.571599999999
$ws.Range("M126").Value = -2365.571599999999
$ws.Range("H132").Value = 4405.1724
$ws.Range("I132").Value = 4707.143
$ws.Range("K132").Value = 14121.429
$ws.Range("M132").Value = -11591.429
$ws.Range("H134").Value = 64111.5
$ws.Range("J134").Value = 64111.5
$ws.Range("L134").Value = 192334.5
$ws.Range("N134").Value = -197404.5
$ws.Range("H136").Value = 26318834
$ws.Range("J136").Value = 5896.923
$ws.Range("L136").Value = 17690.769
$ws.Range("N136").Value = -22790.769
$ws.Range("H138").Value = 99839
$ws.Range("J138").Value = 99839
$ws.Range("L138").Value = 99839
$ws.Range("N138").Value = -110119
